# Insert a new data row at row 114 (pushing existing rows 114..238 down to 115..239)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(114).Insert()

$ws.Range("A114").Value = 4
$ws.Range("B114").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C114").Value = "Los Lagos"
$ws.Range("D114").Value = 44629
$ws.Range("E114").Value = 10
$ws.Range("F114").Value = 100112043
$ws.Range("G114").Value = "Pepino ensalada"
$ws.Range("H114").Value = "Sin especificar"
$ws.Range("I114").Value = "Primera"
$ws.Range("J114").Value = 70
$ws.Range("K114").Value = 22000
$ws.Range("L114").Value = 22000
$ws.Range("M114").Value = 22000
$ws.Range("N114").Value = "$/caja 60 unidades"
$ws.Range("O114").Value = "Región de Arica y Parinacota"
$ws.Range("P114").Value = 367
$ws.Range("Q114").Value = 60
$ws.Range("R114").Value = "Hortaliza"

# Ensure style s="2" (date format) on D114 matches the rest of column D
$ws.Range("D114").NumberFormat = $ws.Range("D115").NumberFormat
